# feat: add 2022-Q3 data
#
# - Insert a brand-new worksheet "2022-Q3" right after "总计", copying the
#   layout/formatting of the existing "2022-Q2" sheet and filling in the
#   new quarter's fund data.
# - Add a corresponding summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

$zongji = $wb.Worksheets.Item("总计")
$q2sheet = $wb.Worksheets.Item("2022-Q2")

# --- 1. Create the new "2022-Q3" sheet right after "总计" -----------------
# Copying the "2022-Q2" sheet preserves every style (header bold+border,
# number formats, page margins, etc.) exactly, so we only need to overwrite
# the cell values afterwards.
$q2sheet.Copy($null, $zongji)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

$newSheet.Range("A2").Value = 0

# Fund code, name, size, position % etc. are stored as plain text in this
# workbook (not numbers) -- force text format while assigning so Excel
# doesn't silently coerce numeric-looking strings ("486002", "3.72", ...)
# into numbers, then restore the default "Normal" style so no stray
# formatting is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $newSheet.Range("B2") "486002"
Set-TextValue $newSheet.Range("C2") "工银全球精选股票（QDII）"
Set-TextValue $newSheet.Range("D2") "3.72"
Set-TextValue $newSheet.Range("E2") "93.69"
Set-TextValue $newSheet.Range("F2") "1.67"
Set-TextValue $newSheet.Range("G2") "0.0621"
$newSheet.Range("H2").Value = 7

# Restore the originally-active sheet (copying changes which tab is
# selected, and the last sheet "2020-Q4" should stay the active one).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

# --- 2. Insert a new row into "总计" for 2022-Q3 ---------------------------
# Shift existing data rows 2:4 down to 3:5 (copying values+formatting),
# then overwrite row 2 with the new quarter's summary figures.
$zongji.Range("A4:D4").Copy($zongji.Range("A5:D5"))
$zongji.Range("A3:D3").Copy($zongji.Range("A4:D4"))
$zongji.Range("A2:D2").Copy($zongji.Range("A3:D3"))

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q3"
$zongji.Range("C2").Value = 1
$zongji.Range("D2").Value = 0.06

$zongji.Range("A3").Value = 1
$zongji.Range("B3").Value = "2022-Q2"
$zongji.Range("C3").Value = 1
$zongji.Range("D3").Value = 0.06

$zongji.Range("A4").Value = 2
$zongji.Range("B4").Value = "2021-Q3"
$zongji.Range("C4").Value = 1
$zongji.Range("D4").Value = 0.08

$zongji.Range("A5").Value = 3
$zongji.Range("B5").Value = "2020-Q4"
$zongji.Range("C5").Value = 1
$zongji.Range("D5").Value = 0.1
